# Update Name of Algo
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = -21.15560000000001
$ws.Range("A10").Value = -20.53969999999997
$ws.Range("A12").Value = -22.49210000000003
$ws.Range("B13").Value = 6.228299999999996
$ws.Range("A18").Value = -22.45920000000003
$ws.Range("A25").Value = -22.33130000000003

$wb.Save()
